# Resilience4J slide: italicize the "Recursos modularizados" bullet list
# (CircuitBreaker, Bulkhead, RateLimiter, Retry, TimeLimiter, Cache).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$targets = @("CircuitBreaker", "Bulkhead", "RateLimiter", "Retry", "TimeLimiter", "Cache")

for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($targets -contains $para.Text.Trim()) {
        $para.Font.Italic = $true
    }
}
